$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 149, shifting existing rows 149:241 down to 150:242.
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new weekly record.
$ws.Cells.Item(149, 1).Value = 5
$ws.Cells.Item(149, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(149, 3).Value = "Maule"
$ws.Cells.Item(149, 4).Value = 44596
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 5).Value = 7
$ws.Cells.Item(149, 6).Value = 100112003
$ws.Cells.Item(149, 7).Value = "Ajo"
$ws.Cells.Item(149, 8).Value = "Chino"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 200
$ws.Cells.Item(149, 11).Value = 20000
$ws.Cells.Item(149, 12).Value = 20000
$ws.Cells.Item(149, 13).Value = 20000
$ws.Cells.Item(149, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(149, 15).Value = "China"
$ws.Cells.Item(149, 16).Value = 2000
$ws.Cells.Item(149, 17).Value = 10
$ws.Cells.Item(149, 18).Value = "Hortaliza"
